$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = " "
$ws.Range("B6").Value = "\31 52174-case-656"
$ws.Range("C6").Value = "Black"

$ws.Range("A7").Value = " "
$ws.Range("B7").Value = "\31 52174-case-656"
$ws.Range("C7").Value = "Black"
